$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.872.22"
$ws.Range("D3").Value = "1.878.79"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.13"
$ws.Range("E5").Value = "  +3.50%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4744"
$ws.Range("E7").Value = "  +5.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3971"
$ws.Range("E8").Value = "  +3.63%  "
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08065"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.89"
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("D13").Value = "1.886.13"
$ws.Range("E13").Value = "  +2.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.965"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001052"
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.26"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06620"
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.28"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "27.952.98"
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.518"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.316"
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("D26").Value = "2.122.59"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.79"
$ws.Range("E27").Value = "  +4.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.28"
$ws.Range("E28").Value = "  +4.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.106"
$ws.Range("E29").Value = "  +2.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.635"
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.61"
$ws.Range("E31").Value = "  +2.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9890"
$ws.Range("E32").Value = "  +6.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09586"
$ws.Range("E33").Value = "  +2.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.469"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.604"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.338"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02275"
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06126"
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.241"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6036"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("E43").Value = "  +3.27%  "
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.276"
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5726"
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.34"
$ws.Range("E47").Value = "  +1.98%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.419"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.954"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06838"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.61"
$ws.Range("E51").Value = "  +5.21%  "
